$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.547.67"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "2.408.26"
$ws.Range("E3").Value = "  +2.90%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'551.85"
$ws.Range("E5").Value = "  +1.87%  "

$ws.Range("D6").Value = "'136.82"
$ws.Range("E6").Value = "  +1.44%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("E9").Value = "  +5.63%  "

$ws.Range("D10").Value = "'5.81"
$ws.Range("E10").Value = "  +5.10%  "

$ws.Range("D11").Value = "'0.361"
$ws.Range("E11").Value = "  +1.33%  "

$ws.Range("E12").Value = "  -2.04%  "

$ws.Range("D13").Value = "'24.70"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").Value = "2.837.61"
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").Value = "59.439.44"
$ws.Range("E15").Value = "  +2.33%  "

$ws.Range("D16").Value = "'0.0000140"
$ws.Range("E16").Value = "  +4.46%  "

$ws.Range("D17").Value = "2.391.93"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("D18").Value = "'11.35"
$ws.Range("E18").Value = "  +6.05%  "

$ws.Range("D19").Value = "'4.44"
$ws.Range("E19").Value = "  +4.23%  "

$ws.Range("D20").Value = "'335.12"

$ws.Range("D21").Value = "'7.01"
$ws.Range("E21").Value = "  +4.78%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "'64.83"
$ws.Range("E23").Value = "  +3.23%  "

$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "  +0.55%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'8.50"
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "'1.37"
$ws.Range("E27").Value = "  -3.38%  "

$ws.Range("D28").Value = "0.0₃0776"
$ws.Range("E28").Value = "  +5.65%  "

$ws.Range("D30").Value = "'170.60"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "'6.26"
$ws.Range("E31").Value = "  +2.58%  "

$ws.Range("D32").Value = "'18.69"
$ws.Range("E32").Value = "  +1.04%  "

$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D35").Value = "'4.30"
$ws.Range("E35").Value = "  +0.57%  "

$ws.Range("D36").Value = "'1.30"
$ws.Range("E36").Value = "  +4.45%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "'1.64"

$ws.Range("D39").Value = "'40.15"
$ws.Range("E39").Value = "  +2.87%  "

$ws.Range("D40").Value = "'0.420"
$ws.Range("E40").Value = "  +11.85%  "

$ws.Range("D41").Value = "'303.25"
$ws.Range("E41").Value = "  +6.48%  "

$ws.Range("E42").Value = "  +2.61%  "

$ws.Range("D43").Value = "'143.01"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0526"
$ws.Range("E44").Value = "  +4.68%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0961"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("D46").Value = "'0.571"
$ws.Range("E46").Value = "  +1.39%  "

$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").Value = "'0.406"
$ws.Range("E47").Value = "  +6.51%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'19.00"
$ws.Range("E48").Value = "  -0.89%  "

$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("D51").Value = "'1.59"
$ws.Range("E51").Value = "  +3.86%  "
